$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8258586525917053
$ws.Range("B1").Value = 5.263033866882324
$ws.Range("C1").Value = 3.676397562026978
$ws.Range("D1").Value = 2.267486095428467
$ws.Range("E1").Value = 1.937081694602966
